function Set-ParaText {
    param($para, $newText)
    $rng = $para.Range
    $rng.End = $rng.End - 1  # exclude the trailing paragraph mark
    $rng.Delete()
    $insRng = $word.ActiveDocument.Range($rng.Start, $rng.Start)
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insRng.InsertXML($xmlFrag)
}

$d = $word.ActiveDocument

# 1) "This is an annotatable resource in the casebook.\n" -> drop the trailing line break
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*This is an annotatable resource in the casebook.*") {
        Set-ParaText $p "This is an annotatable resource in the casebook."
        break
    }
}

# 2) Wrapped multi-line paragraph -> join onto a single line with spaces
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*highlighted: content to highlight*") {
        $newText = "highlighted: content to highlight; elided: content to elide; replaced: content to replace; linked: content to link; noted: content to note; highlighted2: second highlight content;"
        Set-ParaText $p $newText
        break
    }
}

# 3) "This is the second chapter of the casebook.\n" -> drop the trailing line break
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*This is the second chapter of the casebook.*") {
        Set-ParaText $p "This is the second chapter of the casebook."
        break
    }
}
